$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (lambda = 0.02)
$ws.Range("B2").Value = 0.1002355555555555
$ws.Range("C2").Value = 0.000316558558665647
$ws.Range("D2").Value = 900000
$ws.Range("E2").Value = 0.0952617871626035
$ws.Range("F2").Value = 0.0002825813553880523
$ws.Range("G2").Value = 1079331
$ws.Range("H2").Value = 0.01595803326319729
$ws.Range("I2").Value = 0.0001206199565628483
$ws.Range("J2").Value = 1079331
$ws.Range("K2").Value = 0.01132876132421023
$ws.Range("L2").Value = 0.0000813425463776297
$ws.Range("M2").Value = "run_p=0.02_sim=50000.json"

# Row 3 (lambda = 0.1)
$ws.Range("B3").Value = 0.8521133333333333
$ws.Range("C3").Value = 0.0003741897921419155
$ws.Range("D3").Value = 900000
$ws.Range("E3").Value = 0.4019729960584148
$ws.Range("F3").Value = 0.000210977949417813
$ws.Range("G3").Value = 5400619
$ws.Range("H3").Value = 0.05253453354143293
$ws.Range("I3").Value = 0.00009600246074253525
$ws.Range("J3").Value = 5400619
$ws.Range("K3").Value = 0.05824890767162477
$ws.Range("L3").Value = 0.00008615328976067394
$ws.Range("M3").Value = "run_p=0.1_sim=50000.json"

# Row 4 (lambda = 0.2)
$ws.Range("B4").Value = 0.9973588888888889
$ws.Range("C4").Value = 0.00005410006821940513
$ws.Range("D4").Value = 900000
$ws.Range("E4").Value = 0.6463261083130704
$ws.Range("F4").Value = 0.0001454973372395505
$ws.Range("G4").Value = 10798032
$ws.Range("H4").Value = 0.1051356395313516
$ws.Range("I4").Value = 0.00009334289822630437
$ws.Range("J4").Value = 10798032
$ws.Range("K4").Value = 0.1122863190711461
$ws.Range("L4").Value = 0.00008987246824041488
$ws.Range("M4").Value = "run_p=0.2_sim=50000.json"

# Row 5 (lambda = 0.5)
$ws.Range("D5").Value = 900000
$ws.Range("E5").Value = 0.9323585417411794
$ws.Range("F5").Value = 0.0000483370582491058
$ws.Range("G5").Value = 26992026
$ws.Range("H5").Value = 0.3403622240138625
$ws.Range("I5").Value = 0.00009120227967027469
$ws.Range("J5").Value = 26992026
$ws.Range("K5").Value = 0.2702761783178307
$ws.Range("L5").Value = 0.00009577233979532408
$ws.Range("M5").Value = "run_p=0.5_sim=50000.json"

# Row 6 (lambda = 1)
$ws.Range("D6").Value = 900000
$ws.Range("E6").Value = 0.9953703703703703
$ws.Range("F6").Value = 0.000009237800923903298
$ws.Range("G6").Value = 54000000
$ws.Range("H6").Value = 0.575925925925926
$ws.Range("I6").Value = 0.00006725232389930653
$ws.Range("J6").Value = 54000000
$ws.Range("K6").Value = 0.5939393939393953
$ws.Range("L6").Value = 0.0000000000000000004965117957922383
$ws.Range("M6").Value = "run_p=1_sim=50000.json"
